# everyday data update: add a new 2021/11/25 (110年11月25日) row to the top
# of the four daily-series sheets, pushing the existing history down by one
# row (each sheet keeps all of its prior rows, just shifted down by one).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 台指期換倉成本計算 (A1:F7 -> A1:F8) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2:F2").Copy()
$ws1.Range("A2:F2").Insert()
$ws1.Range("A2").Value = "日期：2021/11/25"
$ws1.Range("C2").Value = 17627
$ws1.Range("D2").Value = 5325
$ws1.Range("E2").Value = 3930821
$ws1.Range("F2").Value = 17745

# --- Sheet 2: 散戶多空力道 (A1:B22 -> A1:B23) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2:B2").Copy()
$ws2.Range("A2:B2").Insert()
$ws2.Range("A2").Value = "日期：2021/11/25"
$ws2.Range("B2").Value = 0.21

# --- Sheet 3: 三大法人買賣金額 (A1:C22 -> A1:C23) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2:C2").Copy()
$ws3.Range("A2:C2").Insert()
$ws3.Range("A2").Value = "110年11月25日"
$ws3.Range("B2").Value = -82.64
$ws3.Range("C2").Value = 94.31

# --- Sheet 4: 大盤多空點位 (A1:B21 -> A1:B22) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A2:B2").Copy()
$ws4.Range("A2:B2").Insert()
$ws4.Range("A2").Value = "110年11月25日"
$ws4.Range("B2").Value = 17659.83
